$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.36
$ws.Range("AT2").Value = 2.63
$ws.Range("M3").Value = 1.08
$ws.Range("O3").Value = 1.36
$ws.Range("AT3").Value = 2.63
$ws.Range("J4").Value = 2.38
$ws.Range("M4").Value = 1.08
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.08
$ws.Range("S7").Value = 1.36
$ws.Range("S8").Value = 1.4
$ws.Range("N10").Value = 14.3
$ws.Range("P10").Value = 4.2
$ws.Range("S12").Value = 1.3
$ws.Range("G13").Value = 3.7
$ws.Range("H13").Value = 4.2
$ws.Range("I13").Value = 1.83
$ws.Range("L13").Value = 2.3
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 23
$ws.Range("S13").Value = 1.22
$ws.Range("U13").Value = 1.44
$ws.Range("V13").Value = 2.63
$ws.Range("AC13").Value = 23
$ws.Range("AD13").Value = 9
$ws.Range("AE13").Value = 12
$ws.Range("AI13").Value = 12
$ws.Range("AK13").Value = 17
$ws.Range("AN13").Value = 6.5
$ws.Range("AO13").Value = 19
$ws.Range("AP13").Value = 21
$ws.Range("AX13").Value = 9
$ws.Range("G14").Value = 3.8
$ws.Range("I14").Value = 1.73
$ws.Range("J14").Value = 4
$ws.Range("L14").Value = 2.2
$ws.Range("N14").Value = 23
$ws.Range("S14").Value = 1.2
$ws.Range("AA14").Value = 26
$ws.Range("AX14").Value = 8.5
$ws.Range("I15").Value = 3.6
$ws.Range("L15").Value = 3.75
$ws.Range("M15").Value = 1.02
$ws.Range("N15").Value = 19
$ws.Range("Q15").Value = 1.48
$ws.Range("R15").Value = 2.6
$ws.Range("T15").Value = 3.75
$ws.Range("U15").Value = 1.44
$ws.Range("V15").Value = 2.63
$ws.Range("W15").Value = 12
$ws.Range("X15").Value = 12
$ws.Range("AC15").Value = 19
$ws.Range("AH15").Value = 17
$ws.Range("AP15").Value = 15
$ws.Range("AT15").Value = 3.75
$ws.Range("BA15").Value = 51
$ws.Range("BB15").Value = 101
$ws.Range("G16").Value = 2.1
$ws.Range("I16").Value = 3
$ws.Range("L16").Value = 3.5
$ws.Range("O16").Value = 1.17
$ws.Range("P16").Value = 5
$ws.Range("Q16").Value = 1.57
$ws.Range("R16").Value = 2.35
$ws.Range("S16").Value = 1.29
$ws.Range("T16").Value = 3.5
$ws.Range("X16").Value = 12
$ws.Range("Y16").Value = 9
$ws.Range("AI16").Value = 19
$ws.Range("AK16").Value = 34
$ws.Range("AM16").Value = 26
$ws.Range("AO16").Value = 11
$ws.Range("AP16").Value = 17
$ws.Range("AQ16").Value = 34
$ws.Range("AT16").Value = 3.5
$ws.Range("AW16").Value = 5.5
$ws.Range("AZ16").Value = 51
$ws.Range("BB16").Value = 126
$ws.Range("G17").Value = 1.53
$ws.Range("I18").Value = 2.88
$ws.Range("G20").Value = 1.73
$ws.Range("S21").Value = 1.25
$ws.Range("J22").Value = 1.73
$ws.Range("K22").Value = 2.88
$ws.Range("S22").Value = 1.2
$ws.Range("S23").Value = 1.5
$ws.Range("S25").Value = 1.4
$ws.Range("G27").Value = 1.9
$ws.Range("H27").Value = 3.2
$ws.Range("I27").Value = 4.33
$ws.Range("J27").Value = 2.63
$ws.Range("K27").Value = 2
$ws.Range("Q27").Value = 2.35
$ws.Range("R27").Value = 1.57
$ws.Range("U27").Value = 2.1
$ws.Range("V27").Value = 1.67
$ws.Range("AC27").Value = 7
$ws.Range("AE27").Value = 19
$ws.Range("AG27").Value = 1250
$ws.Range("AM27").Value = 51
$ws.Range("AO27").Value = 11
$ws.Range("AQ27").Value = 41
$ws.Range("BB27").Value = 351
$ws.Range("Q28").Value = 1.98
$ws.Range("R28").Value = 1.88
$ws.Range("AM31").Value = 21
$ws.Range("M32").Value = 1.07
$ws.Range("N32").Value = 9
$ws.Range("Q32").Value = 2.2
$ws.Range("R32").Value = 1.65
$ws.Range("J35").Value = 2.88
$ws.Range("M35").Value = 1.06
$ws.Range("N35").Value = 10
$ws.Range("AE35").Value = 13
$ws.Range("AF35").Value = 41
